try {
    $wb = $excel.ActiveWorkbook
    $ws = $wb.ActiveSheet

    # Row 4 Date Accepted ("1970-01-01") must stay text, not get parsed as a date.
    $ws.Range("H4").NumberFormat = "@"

    # Row 2
    $ws.Range("D2").Value = 'Meteorological parameters are the important factors influencing the infectious diseases such as severe acute respiratory syndrome (SARS) and influenza.
 This study aims to explore the association between Corona Virus Disease 2019 (COVID-19) deaths and weather parameters.
 In this study, we collected the daily death numbers of COVID-19, meteorological parameters and air pollutant data from 20 January 2020 to 29 February 2020 in Wuhan, China.
 Generalized additive model was applied to explore the effect of temperature, humidity and diurnal temperature range on the daily death counts of COVID-19. There were 2299 COVID-19 death counts in Wuhan during the study period.
 A positive association with COVID-19 daily death counts was observed for diurnal temperature range (r = 0.44), but negative association for relative humidity (r = −0.32).
 In addition, one unit increase in diurnal temperature range was only associated with a 2.92% (95% CI: 0.61%, 5.28%) increase in COVID-19 deaths in lag 3. However, both 1 unit increase of temperature and absolute humidity were related to the decreased COVID-19 death in lag 3 and lag 5, with the greatest decrease both in lag 3 [−7.50% (95% CI: −10.99%, −3.88%) and −11.41% (95% CI: −19.68%, −2.29%)].
 In summary, this study suggests the temperature variation and humidity may also be important factors affecting the COVID-19 mortality.
'
    $ws.Range("E2").Value = '[Yueling%Ma%NULL%0, Yadong%Zhao%NULL%1, Jiangtao%Liu%NULL%1, Xiaotao%He%NULL%1, Bo%Wang%NULL%0, Shihua%Fu%NULL%1, Jun%Yan%NULL%1, Jingping%Niu%NULL%1, Ji%Zhou%NULL%1, Bin%Luo%luob@lzu.edu.cn%1]'
    $ws.Range("I2").Value = ''
    $ws.Range("J2").Value = 'Elsevier B.V.'

    # Row 3
    $ws.Range("D3").Value = 'In this study, we aimed at analyzing the associations between transmission of and deaths caused by SARS-CoV-2 and meteorological variables, such as average temperature, minimum temperature, maximum temperature, and precipitation.
 Two outcome measures were considered, with the first aiming to study SARS-CoV-2 infections and the second aiming to study COVID-19 mortality.
 Daily data as well as data on SARS-CoV-2 infections and COVID-19 mortality obtained between December 1, 2019 and March 28, 2020 were collected from weather stations around the world.
 The country''s population density and time of exposure to the disease were used as control variables.
 Finally, a month dummy variable was added.
 Daily data by country were analyzed using the panel data model.
 An increase in the average daily temperature by one degree Fahrenheit reduced the number of cases by approximately 6.4 cases/day.
 There was a negative correlation between the average temperature per country and the number of cases of SARS-CoV-2 infections.
 This association remained strong even with the incorporation of additional variables and controls (maximum temperature, average temperature, minimum temperature, and precipitation) and fixed country effects.
 There was a positive correlation between precipitation and SARS-CoV-2 transmission.
 Countries with higher rainfall measurements showed an increase in disease transmission.
 For each average inch/day, there was an increase of 56.01 cases/day.
 COVID-19 mortality showed no significant association with temperature.
'
    $ws.Range("E3").Value = '[Marcos Felipe Falcão%Sobral%NULL%0, Gisleia Benini%Duarte%NULL%1, Ana Iza Gomes%da Penha Sobral%NULL%1, Marcelo Luiz Monteiro%Marinho%NULL%1, André%de Souza Melo%NULL%1]'
    $ws.Range("I3").Value = ''
    $ws.Range("J3").Value = 'Elsevier B.V.'

    # Row 4
    $ws.Range("C4").Value = 'Influence of socio-ecological factors on COVID-19 risk: a cross-sectional study based on 178 countries/regions worldwide'
    $ws.Range("D4").Value = 'Background:
id="P3">The initial outbreak of COVID-19 caused by SARS-CoV-2 in China in 2019 has been severely tested in other countries worldwide.
 We aimed to describe the spatial distribution of the COVID-19 pandemic worldwide and assess the effects of various socio-ecological factors on COVID-19 risk.
Methods:
id="P4">We collected COVID-19 pandemic infection data and social-ecological data of 178 countries/regions worldwide from three database.
 We used spatial econometrics method to assess the global and local correlation of COVID-19 risk indicators for COVID-19. To estimate the adjusted incidence rate ratio (IRR), we modelled negative binomial regression analysis with spatial information and socio-ecological factors.
Findings:
id="P5">The study indicated that 37, 29 and 39 countries/regions were strongly opposite from the IR, CMR and DCI index "spatial autocorrelation hypothesis", respectively.
 The IRs were significantly positively associated with GDP per capita, the use of at least basic sanitation services and social insurance program coverage, and were significantly negatively associated with the proportion of the population spending more than 25% of household consumption or income on out-of-pocket health care expenses and the poverty headcount ratio at the national poverty lines.
 The CMR was significantly positively associated with urban populations, GDP per capita and current health expenditure, and was significantly negatively associated with the number of hospital beds, number of nurses and midwives, and poverty headcount ratio at the national poverty lines.
 The DCI was significantly positively associated with urban populations, population density and researchers in R&amp;D, and was significantly negatively associated with the number of hospital beds, number of nurses and midwives and poverty headcount ratio at the national poverty lines.
 We also found that climatic factors were not significantly associated with COVID-19 risk.
Conclusion:
id="P6">Countries/regions should pay more attention to controlling population flow, improving diagnosis and treatment capacity, and improving public welfare policies.
'
    $ws.Range("E4").Value = '[Dai%Su%NULL%0, Yingchun%Chen%NULL%1, Kevin%He%NULL%1, Tao%Zhang%NULL%0, Min%Tan%NULL%1, Yunfan%Zhang%NULL%1, Xingyu%Zhang%NULL%1]'
    $ws.Range("F4").Value = 'PMC7276015'
    $ws.Range("G4").Value = 'PMC'
    $ws.Range("H4").Value = '1970-01-01'
    $ws.Range("I4").Value = ''
    $ws.Range("J4").Value = 'Cold Spring Harbor Laboratory'

    # Row 5
    $ws.Range("D5").Value = 'The coronavirus disease 2019 (COVID-19) pandemic is the defining global health crisis of our time and the greatest challenge facing the world.
 Meteorological parameters are reportedly crucial factors affecting respiratory infectious disease epidemics; however, the effect of meteorological parameters on COVID-19 remains controversial.
 This study investigated the effects of temperature and relative humidity on daily new cases and daily new deaths of COVID-19, which has useful implications for policymakers and the public.
 Daily data on meteorological conditions, new cases and new deaths of COVID-19 were collected for 166 countries (excluding China) as of March 27, 2020. Log-linear generalized additive model was used to analyze the effects of temperature and relative humidity on daily new cases and daily new deaths of COVID-19, with potential confounders controlled for, including wind speed, median age of the national population, Global Health Security Index, Human Development Index and population density.
 Our findings revealed that temperature and relative humidity were both negatively related to daily new cases and deaths.
 A 1 °C increase in temperature was associated with a 3.08% (95% CI: 1.53%, 4.63%) reduction in daily new cases and a 1.19% (95% CI: 0.44%, 1.95%) reduction in daily new deaths, whereas a 1% increase in relative humidity was associated with a 0.85% (95% CI: 0.51%, 1.19%) reduction in daily new cases and a 0.51% (95% CI: 0.34%, 0.67%) reduction in daily new deaths.
 The results remained robust when different lag structures and the sensitivity analysis were used.
 These findings provide preliminary evidence that the COVID-19 pandemic may be partially suppressed with temperature and humidity increases.
 However, active measures must be taken to control the source of infection, block transmission and prevent further spread of COVID-19.'
    $ws.Range("E5").Value = '[Yu%Wu%NULL%0, Wenzhan%Jing%NULL%1, Jue%Liu%NULL%1, Qiuyue%Ma%NULL%1, Jie%Yuan%NULL%1, Yaping%Wang%NULL%1, Min%Du%NULL%1, Min%Liu%NULL%0]'
    $ws.Range("I5").Value = ''
    $ws.Range("J5").Value = 'The Authors. Published by Elsevier B.V.'

    # Row 6
    $ws.Range("E6").Value = '[Yasir%Rehman%NULL%0, Nadia%Rehman%NULL%1]'
    $ws.Range("I6").Value = ''
    $ws.Range("J6").Value = 'AIMS Press'

    # Row 7
    $ws.Range("D7").Value = 'Novel corona virus disease 2019 (COVID-19), which first emerged in December 2019, has become a pandemic.
 This study aimed to investigate the associations between meteorological factors and COVID-19 incidence and mortality worldwide.
 This study included 1,908,197 confirmed cases of and 119,257 deaths from COVID-19 from 190 countries between 23 January and 13 April, 2020. We used a distributed lag non-linear model with city-/country-level random intercept to investigate the associations between COVID19 incidence and daily temperature, relative humidity, and wind speed.
 A series of confounders were considered in the analysis including demographics, socioeconomics, geographic locations, and political strategies.
 Sensitivity analyses were performed to examine the robustness of the associations.
 The COVID-19 incidence showed a stronger association with temperature than with relative humidity or wind speed.
 An inverse association was identified between the COVID-19 incidence and temperature.
 The corresponding 14-day cumulative relative risk was 1.28 [95% confidence interval (CI), 1.20–1.36] at 5 °C, and 0.75 (95% CI, 0.65–0.86) at 22 °C with reference to the risk at 11 °C.
 An inverse J-shaped association was observed between relative humidity and the COVID-19 incidence, with the highest risk at 72%.
 A higher wind speed was associated with a generally lower incidence of COVID-19, although the associations were weak.
 Sensitivity analyses generally yielded similar results.
 The COVID-19 incidence decreased with the increase of temperature.
 Our study suggests that the spread of COVID-19 may slow during summer but may increase during winter.
'
    $ws.Range("E7").Value = '[Cui%Guo%NULL%0, Yacong%Bo%NULL%1, Changqing%Lin%NULL%1, Hao Bi%Li%NULL%1, Yiqian%Zeng%NULL%1, Yumiao%Zhang%NULL%1, Md Shakhaoat%Hossain%NULL%1, Jimmy W.M.%Chan%NULL%1, David W.%Yeung%NULL%1, Kin-on%Kwok%NULL%1, Samuel Y.S.%Wong%NULL%1, Alexis K.H.%Lau%NULL%1, Xiang Qian%Lao%NULL%1]'
    $ws.Range("I7").Value = ''
    $ws.Range("J7").Value = 'Elsevier B.V.'

    # Row 8
    $ws.Range("D8").Value = 'id="Par1">Novel coronavirus (SARS-CoV-2) causing COVID-19 disease has arisen to be a pandemic.
 Since there is a close association between other viral infection cases by epidemics and environmental factors, this study intends to unveil meteorological effects on the outbreak of COVID-19 across eight divisions of Bangladesh from March to April 2020. A compound Poisson generalized linear modeling (CPGLM), along with a Monte-Carlo method and random forest (RF) model, was employed to explore how meteorological factors affecting the COVID-19 transmission in Bangladesh.
 Results showed that subtropical climate (mean temperature about 26.6 °C, mean relative humidity (MRH) 64%, and rainfall approximately 3 mm) enhanced COVD-19 onset.
 The CPGLM model revealed that every 1 mm increase in rainfall elevated by 30.99% (95% CI 77.18%, − 15.20%) COVID-19 cases, while an increase of 1 °C of diurnal temperature (TDN) declined the confirmed cases by − 14.2% (95% CI 9.73%, − 38.13%) on the lag 1 and lag 2, respectively.
 In addition, NRH and MRH had the highest increase (17.98% (95% CI 22.5%, 13.42%) and 19.92% (95% CI: 25.71%, 14.13%)) of COVID-19 cased in lag 4. The results of the RF model indicated that TDN and AH (absolute humidity) influence the COVID-19 cases most.
 In the Dhaka division, MRH is the most vital meteorological factor that affects COVID-19 deaths.
 This study indicates the humidity and rainfall are crucial factors affecting the COVID-19 case, which is contrary to many previous studies in other countries.
 These outcomes can have policy formulation for the suppression of the COVID-19 outbreak in Bangladesh.
'
    $ws.Range("E8").Value = '[Abu Reza Md. Towfiqul%Islam%towfiq_dm@brur.ac.bd%0, Md.%Hasanuzzaman%NULL%2, Md.%Hasanuzzaman%NULL%0, Mashura%Shammi%NULL%1, Roquia%Salam%NULL%1, Md.%Bodrud-Doza%NULL%1, Md. Mostafizur%Rahman%rahmanmm@juniv.edu%1, Md. Abdul%Mannan%NULL%1, Saleemul%Huq%NULL%1]'
    $ws.Range("I8").Value = ''
    $ws.Range("J8").Value = 'Springer Berlin Heidelberg'

    # Row 9
    $ws.Range("D9").Value = 'id="Par1">The emergence of coronavirus disease 2019 (COVID-19) has become a worldwide pandemic after its first outbreak in Wuhan, China.
 However, it remains unclear whether COVID-19 death is linked to ambient air pollutants or meteorological conditions.
 We collected the daily COVID-19 death number, air quality index (AQI), ambient air pollutant concentrations, and meteorological variables data of Wuhan between Jan 25 and April 7, 2020. The Pearson and Poisson regression models were used accordingly to understand the association between COVID-19 deaths and each risk factor.
 The daily COVID-19 deaths were positively correlated with AQI (slope = 0.4 ± 0.09, R2 = 0.24, p &lt; 0.01).
 Detailedly, PM2.5 was the only pollutant exhibiting a positive association (relative risk (RR) = 1.079, 95%CI 1.071–1.086, p &lt; 0.01) with COVID-19 deaths.
 The PM10, SO2, and CO were all also significantly associated with COVID-19 deaths, but in negative pattern (p &lt; 0.01).
 Among them, PM10 and CO had the highest and lowest RR, which equaled to 0.952 (95%CI 0.945–0.959) and 0.177 (95%CI 0.131–0.24), respectively.
 Additionally, temperature was inversely associated with COVID-19 deaths (RR = 0.861, 95%CI 0.851–0.872, p &lt; 0.01).
 Contrarily, diurnal temperature range was positively associated with COVID-19 deaths (RR = 1.014, 95%CI 1.003–1.025, p &lt; 0.05).
 The data suggested that PM2.5 and diurnal temperature range are tightly associated with increased COVID-19 deaths.
'
    $ws.Range("E9").Value = '[Ying%Jiang%yjiang1@vt.edu%0, Jing%Xu%xujingjingalpha@sohu.com%0]'
    $ws.Range("I9").Value = ''
    $ws.Range("J9").Value = 'Springer Netherlands'

    # Row 10
    $ws.Range("D10").Value = 'In this study, we aimed to examine spatial inequalities of COVID-19 mortality rate in relation to spatial inequalities of socioeconomic and environmental factors across England.
 Specifically, we first explored spatial patterns of COVID-19 mortality rate in comparison to non-COVID-19 mortality rate.
 Subsequently, we established models to investigate contributions of socioeconomic and environmental factors to spatial variations of COVID-19 mortality rate across England (N = 317).
 Two newly developed specifications of spatial regression models were established successfully to estimate COVID-19 mortality rate (R2 = 0.49 and R2 = 0.793).
 The level of spatial inequalities of COVID-19 mortality is higher than that of non-COVID-19 mortality in England.
 Although global spatial association of COVID-19 mortality and non-COVID-19 mortality is positive, local spatial association of COVID-19 mortality and non-COVID-19 mortality is negative in some areas.
 Expectedly, hospital accessibility is negatively related to COVID-19 mortality rate.
 Percent of Asians, percent of Blacks, and unemployment rate are positively related to COVID-19 mortality rate.
 More importantly, relative humidity is negatively related to COVID-19 mortality rate.
 Moreover, among the spatial models estimated, the ‘random effects specification of eigenvector spatial filtering model’ outperforms the ‘matrix exponential spatial specification of spatial autoregressive model’.
'
    $ws.Range("E10").Value = '[Yeran%Sun%NULL%0, Xuke%Hu%NULL%1, Jing%Xie%NULL%0]'
    $ws.Range("I10").Value = ''
    $ws.Range("J10").Value = 'Elsevier B.V.'

    # Row 11
    $ws.Range("D11").Value = 'The transmission of the severe acute respiratory syndrome coronavirus 2 (SARS-CoV-2) and the severity of the related disease (COVID-19) are influenced by a large number of factors.
 This study aimed to investigate the correlation of COVID-19 case and death rates with possible causal climatological and sociodemographic factors for the March to May 2020 (first wave) period in a worldwide scale by statistically processing data for over one hundred countries.
 The weather parameters considered herein were air temperature, relative humidity, cumulative precipitation, and cloud cover, while sociodemographic factors included population density, median age, and government measures in response to the pandemic.
 The results of this study indicate that there is a statistically significant correlation between average atmospheric temperature and the COVID-19 case and death rates, with chi-square test p-values in the 0.001–0.02 range.
 Regarding sociodemographic factors, there is an even stronger dependence of the case and death rates on the population median age (p = 0.0006–0.0012).
 Multivariate linear regression analysis using Lasso and the forward stepwise approach revealed that the median age ranks first in importance among the examined variables, followed by the temperature and the delays in taking first governmental measures or issuing stay-at-home orders.
'
    $ws.Range("E11").Value = '[Ploutarchos%Tzampoglou%NULL%0, Dimitrios%Loukidis%NULL%1]'
    $ws.Range("I11").Value = ''
    $ws.Range("J11").Value = 'MDPI'

    # Row 12
    $ws.Range("D12").Value = 'On March 12th, 2020, the WHO declared COVID-19 as a pandemic.
 The collective impact of environmental and ecosystem factors, as well as biodiversity, on the spread of COVID-19 and its mortality evolution remain empirically unknown, particularly in regions with a wide ecosystem range.
 The aim of our study is to assess how those factors impact on the COVID-19 spread and mortality by country.
 This study compiled a global database merging WHO daily case reports with other publicly available measures from January 21st to May 18th, 2020. We applied spatio-temporal models to identify the influence of biodiversity, temperature, and precipitation and fitted generalized linear mixed models to identify the effects of environmental variables.
 Additionally, we used count time series to characterize the association between COVID-19 spread and air quality factors.
 All analyses were adjusted by social demographic, country-income level, and government policy intervention confounders, among 160 countries, globally.
 Our results reveal a statistically meaningful association between COVID-19 infection and several factors of interest at country and city levels such as the national biodiversity index, air quality, and pollutants elements (PM10, PM2.5, and O3).
 Particularly, there is a significant relationship of loss of biodiversity, high level of air pollutants, and diminished air quality with COVID-19 infection spread and mortality.
 Our findings provide an empirical foundation for future studies on the relationship between air quality variables, a country’s biodiversity, and COVID-19 transmission and mortality.
 The relationships measured in this study can be valuable when governments plan environmental and health policies, as alternative strategy to respond to new COVID-19 outbreaks and prevent future crises.
'
    $ws.Range("E12").Value = '[Daniel%Fernández%NULL%0, Iago%Giné-Vázquez%NULL%1, Ivy%Liu%NULL%1, Recai%Yucel%NULL%1, Marta%Nai Ruscone%NULL%1, Marianthi%Morena%NULL%1, Víctor Gerardo%García%NULL%1, Josep Maria%Haro%NULL%1, William%Pan%NULL%1, Stefanos%Tyrovolas%NULL%1]'
    $ws.Range("I12").Value = ''
    $ws.Range("J12").Value = 'Elsevier Ltd.'

    Write-Host "Citation graph abstracts/authors/publisher fields updated."
} catch {
    Write-Host "ERROR: $_"
    throw
}
